$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.889.33'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '2.255.22'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '2.599.00'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('D15').Value = '2.350.19'
$ws.Range('E15').Value = '  +4.84%  '
$ws.Range('E16').Value = '  +2.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').Value = '44.615.96'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('D19').Value = '0.0₃0958'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  +3.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.88'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.38%  '
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.03'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0801'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').Value = '  +5.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.15'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D44').Value = '1.845.02'
$ws.Range('E44').Value = '  +5.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +18.64%  '
$ws.Range('E46').Value = '  +2.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '80.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '99.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.39%  '
